# Updates per-row currentAveragePrice/currentAveragePriceNQ/currentAveragePriceHQ,
# LevePriceNQ/LevePriceHQ/LeveProfitNQ/LeveProfitHQ figures (columns H-N) pulled by
# the scheduled market-data refresh, sheet by sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1450
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1450
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 4350
$ws.Range("N121").Value = -7844
$ws.Range("M121").Value = $null  # column no longer populated for this row

$ws.Range("H126").Value = 45123.332
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 45123.332
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 45123.332
$ws.Range("N126").Value = -55003.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 37899.332
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 37899.332
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 37899.332
$ws.Range("N44").Value = -38875.332

$ws.Range("H55").Value = 20715.75
$ws.Range("I55").Value = 48
$ws.Range("J55").Value = 23668.285
$ws.Range("K55").Value = 48
$ws.Range("L55").Value = 23668.285
$ws.Range("M55").Value = 267
$ws.Range("N55").Value = -24298.285

$ws.Range("H80").Value = 26654.889
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 26654.889
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 26654.889
$ws.Range("N80").Value = -28650.889

$ws.Range("H83").Value = 26654.889
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 26654.889
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 79964.667
$ws.Range("N83").Value = -89948.667

$ws.Range("H128").Value = 52495
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 52495
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 52495
$ws.Range("N128").Value = -62455

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 34832
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 34832
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 34832
$ws.Range("N35").Value = -35452

$ws.Range("H122").Value = 49437.777
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 49437.777
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 49437.777
$ws.Range("N122").Value = -59237.777

$ws.Range("H125").Value = 47247.5
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 47247.5
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 47247.5
$ws.Range("N125").Value = -57087.5

$ws.Range("H126").Value = 47400
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 47400
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 47400
$ws.Range("N126").Value = -57280

$ws.Range("H132").Value = 37455.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 37455.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 37455.5
$ws.Range("N132").Value = -47575.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49294.363
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 49294.363
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 49294.363
$ws.Range("N20").Value = -49766.363

$ws.Range("H30").Value = 49294.363
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 49294.363
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 49294.363
$ws.Range("N30").Value = -49476.363

$ws.Range("H41").Value = 14569.857
$ws.Range("I41").Value = 5479.5
$ws.Range("J41").Value = 18206
$ws.Range("K41").Value = 5479.5
$ws.Range("L41").Value = 18206
$ws.Range("M41").Value = -5051.5
$ws.Range("N41").Value = -19062

$ws.Range("H68").Value = 34230.6
$ws.Range("I68").Value = 30268
$ws.Range("J68").Value = 35221.25
$ws.Range("K68").Value = 30268
$ws.Range("L68").Value = 35221.25
$ws.Range("M68").Value = -29519
$ws.Range("N68").Value = -36719.25

$ws.Range("H71").Value = 34230.6
$ws.Range("I71").Value = 30268
$ws.Range("J71").Value = 35221.25
$ws.Range("K71").Value = 90804
$ws.Range("L71").Value = 105663.75
$ws.Range("M71").Value = -87060
$ws.Range("N71").Value = -113151.75

$ws.Range("H128").Value = 49294.363
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 49294.363
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 49294.363
$ws.Range("N128").Value = -59254.363

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 6071.684
$ws.Range("I16").Value = 5990
$ws.Range("J16").Value = 6076.222
$ws.Range("K16").Value = 17970
$ws.Range("L16").Value = 18228.666
$ws.Range("M16").Value = -17797
$ws.Range("N16").Value = -18574.666

$ws.Range("H20").Value = 12244
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 12244
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 36732
$ws.Range("N20").Value = -37186

$ws.Range("H22").Value = 1333.3334
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 9000
$ws.Range("M22").Value = -1331
$ws.Range("N22").Value = -9338

$ws.Range("H27").Value = 1333.3334
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 9000
$ws.Range("M27").Value = -1398
$ws.Range("N27").Value = -9204

$ws.Range("H113").Value = 581.17145
$ws.Range("I113").Value = 561.6429000000001
$ws.Range("J113").Value = 594.1905
$ws.Range("K113").Value = 1684.9287
$ws.Range("L113").Value = 1782.5715
$ws.Range("M113").Value = 485.0712999999998
$ws.Range("N113").Value = -6122.5715

$ws.Range("H117").Value = 2972.5757
$ws.Range("I117").Value = 516.6667
$ws.Range("J117").Value = 3518.3333
$ws.Range("K117").Value = 1550.0001
$ws.Range("L117").Value = 10554.9999
$ws.Range("M117").Value = 1891.9999
$ws.Range("N117").Value = -17438.9999

$ws.Range("H129").Value = 1838.5
$ws.Range("I129").Value = 945.55554
$ws.Range("J129").Value = 2569.0908
$ws.Range("K129").Value = 2836.66662
$ws.Range("L129").Value = 7707.2724
$ws.Range("M129").Value = 2163.33338
$ws.Range("N129").Value = -17707.2724

$ws.Range("H132").Value = 823536.4
$ws.Range("I132").Value = 1197100.1
$ws.Range("J132").Value = 1696
$ws.Range("K132").Value = 10773900.9
$ws.Range("L132").Value = 15264
$ws.Range("M132").Value = -10771370.9
$ws.Range("N132").Value = -20324

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 17007.908
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 17007.908
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 17007.908
$ws.Range("N57").Value = -18647.908

$ws.Range("H70").Value = 29369.117
$ws.Range("I70").Value = 35593.293
$ws.Range("J70").Value = 5855.5557
$ws.Range("K70").Value = 35593.293
$ws.Range("L70").Value = 5855.5557
$ws.Range("M70").Value = -35323.293
$ws.Range("N70").Value = -6395.5557

$ws.Range("H73").Value = 29369.117
$ws.Range("I73").Value = 35593.293
$ws.Range("J73").Value = 5855.5557
$ws.Range("K73").Value = 35593.293
$ws.Range("L73").Value = 5855.5557
$ws.Range("M73").Value = -34657.293
$ws.Range("N73").Value = -7727.5557

$ws.Range("H122").Value = 2277.4211
$ws.Range("I122").Value = 2238.2666
$ws.Range("J122").Value = 2424.25
$ws.Range("K122").Value = 6714.7998
$ws.Range("L122").Value = 7272.75
$ws.Range("M122").Value = -4264.7998
$ws.Range("N122").Value = -12172.75

$ws.Range("H127").Value = 48068
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 48068
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 48068
$ws.Range("N127").Value = -57988

$ws.Range("H130").Value = 53570
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 53570
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 53570
$ws.Range("N130").Value = -63610

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 38392.332
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 38392.332
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 38392.332
$ws.Range("N92").Value = -43384.332

$ws.Range("H108").Value = 42542
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 42542
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 42542
$ws.Range("N108").Value = -50222

$ws.Range("H109").Value = 28316.666
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 28316.666
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 28316.666
$ws.Range("N109").Value = -31090.666

$ws.Range("H127").Value = 42336.43
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 42336.43
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 42336.43
$ws.Range("N127").Value = -52256.43

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 29688.5
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 29688.5
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 29688.5
$ws.Range("N109").Value = -32462.5
